$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 256.02563
$ws.Range("I33").Value = 249.85294
$ws.Range("K33").Value = 249.85294
$ws.Range("M33").Value = -20.85293999999999
$ws.Range("H100").Value = 1452.6316
$ws.Range("I100").Value = 1133.3334
$ws.Range("K100").Value = 1133.3334
$ws.Range("M100").Value = -592.3334
$ws.Range("H137").Value = 4142.921
$ws.Range("J137").Value = 2324.0833
$ws.Range("L137").Value = 6972.249899999999
$ws.Range("N137").Value = -12072.2499
# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 49496.617
$ws.Range("I45").Value = 92430.45
$ws.Range("J45").Value = 2269.4
$ws.Range("K45").Value = 92430.45
$ws.Range("L45").Value = 2269.4
$ws.Range("M45").Value = -92053.45
$ws.Range("N45").Value = -3023.4
$ws.Range("H110").Value = 37234.863
$ws.Range("I110").Value = 892.61536
$ws.Range("J110").Value = 89729.22
$ws.Range("K110").Value = 892.61536
$ws.Range("L110").Value = 89729.22
$ws.Range("M110").Value = 1152.38464
$ws.Range("N110").Value = -93819.22
$ws.Range("H122").Value = 2749.9167
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 2833.2222
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 8499.6666
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -13399.6666
$ws.Range("H132").Value = 9955702
$ws.Range("I132").Value = 15581558
$ws.Range("J132").Value = 2264.3076
$ws.Range("K132").Value = 46744674
$ws.Range("L132").Value = 6792.9228
$ws.Range("M132").Value = -46742144
$ws.Range("N132").Value = -11852.9228
# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1464.0526
$ws.Range("I94").Value = 1293.8667
$ws.Range("K94").Value = 1293.8667
$ws.Range("M94").Value = -842.8667
$ws.Range("H134").Value = 49180.87
$ws.Range("I134").Value = 61881.11
$ws.Range("J134").Value = 3460
$ws.Range("K134").Value = 185643.33
$ws.Range("L134").Value = 10380
$ws.Range("M134").Value = -183108.33
$ws.Range("N134").Value = -15450
# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3473850.8
$ws.Range("I31").Value = 1340.762
$ws.Range("J31").Value = 6174692
$ws.Range("K31").Value = 1340.762
$ws.Range("L31").Value = 6174692
$ws.Range("M31").Value = -1045.762
$ws.Range("N31").Value = -6175282
$ws.Range("H34").Value = 3473850.8
$ws.Range("I34").Value = 1340.762
$ws.Range("J34").Value = 6174692
$ws.Range("K34").Value = 1340.762
$ws.Range("L34").Value = 6174692
$ws.Range("M34").Value = -1138.762
$ws.Range("N34").Value = -6175096
$ws.Range("H94").Value = 781.1429000000001
$ws.Range("I94").Value = 1606
$ws.Range("J94").Value = 643.6667
$ws.Range("K94").Value = 1606
$ws.Range("L94").Value = 643.6667
$ws.Range("M94").Value = -1155
$ws.Range("N94").Value = -1545.6667
$ws.Range("H132").Value = 3477.8147
$ws.Range("I132").Value = 3556.5
$ws.Range("J132").Value = 3363.3635
$ws.Range("K132").Value = 10669.5
$ws.Range("L132").Value = 10090.0905
$ws.Range("M132").Value = -8139.5
$ws.Range("N132").Value = -15150.0905
$ws.Range("H134").Value = 2012.3948
$ws.Range("I134").Value = 1956.5518
$ws.Range("K134").Value = 5869.6554
$ws.Range("M134").Value = -3334.6554
# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 801301.75
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 1602599
$ws.Range("K2").Value = 27
$ws.Range("L2").Value = 9615594
$ws.Range("M2").Value = 86
$ws.Range("N2").Value = -9615820
$ws.Range("H17").Value = 294.45456
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 332.1111
$ws.Range("K17").Value = 375
$ws.Range("L17").Value = 996.3333
$ws.Range("M17").Value = -206
$ws.Range("N17").Value = -1334.3333
$ws.Range("H60").Value = 27779546
$ws.Range("I60").Value = 55556010
$ws.Range("J60").Value = 3083.3333
$ws.Range("K60").Value = 166668030
$ws.Range("L60").Value = 9249.999899999999
$ws.Range("M60").Value = -166667779
$ws.Range("N60").Value = -9751.999899999999
$ws.Range("H81").Value = 200001000
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 250001120
$ws.Range("K81").Value = 1500
$ws.Range("L81").Value = 750003360
$ws.Range("M81").Value = -377
$ws.Range("N81").Value = -750005606
$ws.Range("H84").Value = 200001000
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 250001120
$ws.Range("K84").Value = 4500
$ws.Range("L84").Value = 2250010080
$ws.Range("M84").Value = 1116
$ws.Range("N84").Value = -2250021312
$ws.Range("H98").Value = 1180.6
$ws.Range("I98").Value = 1225.75
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 3677.25
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -2179.25
$ws.Range("N98").Value = -5996
$ws.Range("H113").Value = 508.41177
$ws.Range("I113").Value = 701
$ws.Range("J113").Value = 467.14285
$ws.Range("K113").Value = 2103
$ws.Range("L113").Value = 1401.42855
$ws.Range("M113").Value = 67
$ws.Range("N113").Value = -5741.428550000001
$ws.Range("H131").Value = 1668679.4
$ws.Range("J131").Value = 2002139.2
$ws.Range("L131").Value = 6006417.6
$ws.Range("N131").Value = -6016497.6
$ws.Range("H133").Value = 7334.3687
$ws.Range("J133").Value = 7617.2354
$ws.Range("L133").Value = 22851.7062
$ws.Range("N133").Value = -32971.7062
$ws.Range("H140").Value = 25000996
$ws.Range("I140").Value = 25000996
$ws.Range("K140").Value = 75002988
$ws.Range("M140").Value = -74997808
# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1237.5
$ws.Range("I97").Value = 1259.091
$ws.Range("K97").Value = 1259.091
$ws.Range("M97").Value = -763.0909999999999
$ws.Range("H122").Value = 7207.643
$ws.Range("I122").Value = 8090.7
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 24272.1
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -21822.1
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 1409.091
$ws.Range("I126").Value = 1083.3334
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 3250.0002
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -780.0001999999999
$ws.Range("N126").Value = -10340
# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1878.55
$ws.Range("I16").Value = 1878.55
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1878.55
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1708.55
$ws.Range("H40").Value = 3243.348
$ws.Range("I40").Value = 3171.2856
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 3171.2856
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -3035.2856
$ws.Range("N40").Value = -4272
# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1644.3684
$ws.Range("I126").Value = 1941.6428
$ws.Range("J126").Value = 812
$ws.Range("K126").Value = 5824.928400000001
$ws.Range("L126").Value = 2436
$ws.Range("M126").Value = -3354.928400000001
$ws.Range("N126").Value = -7376
$ws.Range("H132").Value = 6765363.5
$ws.Range("I132").Value = 8174446
$ws.Range("J132").Value = 1768.4
$ws.Range("K132").Value = 24523338
$ws.Range("L132").Value = 5305.200000000001
$ws.Range("M132").Value = -24520808
$ws.Range("N132").Value = -10365.2
$ws.Range("H135").Value = 91905
$ws.Range("J135").Value = 91905
$ws.Range("L135").Value = 91905
$ws.Range("N135").Value = -102045
